# Update countries & provincias Spain
# Applies the Oct 23 2020 14:39 data refresh to the "Pais" sheet:
#   - updates case-count figures for several countries
#   - some countries overtake their neighbours in the (descending) sort
#     order, so their rows swap places
#   - refreshes the "last updated" timestamp in the title cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp (row 1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 14:39"

# --- Helper: write a full data row (country + 7 numeric columns) -----
# Columns: A Pais, B Casos totales, C Nuevos casos, D Casos activos,
#          E Recuperados, F Casos criticos, G Muertes hoy, H Muertes
function Set-CountryRow($Row, $Country, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 43: Kuwait - updated figures
Set-CountryRow 43 "Kuwait" 120232 812 111440 8052 0 10 740

# Row 44 (Oman) is unchanged.

# Rows 45-47: Suecia jumps ahead of Kazajistan & Portugal with fresh data;
# Kazajistan and Portugal shift down one row each, keeping their own figures.
Set-CountryRow 45 "Suecia" 110594 0 0 0 0 7 5933
Set-CountryRow 46 "Kazajistan" 110086 179 105493 2797 0 0 1796
Set-CountryRow 47 "Portugal" 109541 0 64531 42765 0 0 2245

# Row 48 (Egipto) is unchanged.

# Row 49: Suiza - updated figures
Set-CountryRow 49 "Suiza" 103653 6634 55800 45796 0 5 2057

# Row 73: Ghana - updated figures
Set-CountryRow 73 "Ghana" 47601 63 46824 463 0 2 314

# Rows 80-81: Dinamarca overtakes Eslovaquia with fresh data;
# Eslovaquia shifts down one row, keeping its own figures.
Set-CountryRow 80 "Dinamarca" 38622 859 31295 6630 0 3 697
Set-CountryRow 81 "Eslovaquia" 37911 2581 8859 28918 0 19 134

# Rows 95-96: Eslovenia overtakes Albania with fresh data;
# Albania shifts down one row, keeping its own figures.
Set-CountryRow 95 "Eslovenia" 19307 1656 7659 11434 0 3 214
Set-CountryRow 96 "Albania" 18250 0 10395 7390 0 0 465

# Rows 142-143: Islandia overtakes Estonia; both get refreshed figures.
Set-CountryRow 142 "Islandia" 4308 40 3187 1110 0 0 11
Set-CountryRow 143 "Estonia" 4300 53 3418 809 0 2 73

# Row 163: Lesoto - updated figures (stays in place)
Set-CountryRow 163 "Lesoto" 1934 11 961 930 0 0 43

# Rows 216-217: Montserrat and Islas Malvinas swap places (no data change).
Set-CountryRow 216 "Montserrat" 13 0 12 0 0 0 1
Set-CountryRow 217 "Islas Malvinas" 13 0 13 0 0 0 0
